# Append two new data rows (6 and 7) to the "Statistic" sheet, following the
# same layout as the existing rows (row 4 = config #2, row 5 = config #3):
#   row 6 -> config #4 (only the first "Прогін 0" block of stats is filled in,
#            plus a leading 0 for the next block's NFE column)
#   row 7 -> config #5 (identical shape, same values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ---------------------------------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 77200
$ws.Range("C6").Value = 413
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 1
$ws.Range("AC6").Value = 0

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 77200
$ws.Range("C7").Value = 413
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 1
$ws.Range("AC7").Value = 0
